# Knowledge Enhancement Plan - update by Kamna on 24 Aug
# Adds a new "Latest status as on 5/24" column (I) with status updates
# for the team, mirrors some existing 5/23 statuses into it, and records
# a couple of new status notes ("On Leave" / GitHub doc update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: same width as column H -----------------------------
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# --- Header (row 1): new "Latest status as on 5/24" header in I1 ------
# Copy H1's look (border/fill/bold header style) into I1, then set text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "Latest status as on 5/24"

# --- Row 2 & 3: carry the 5/23 status forward into the new column -----
# (identical text/format to column H for these two rows)
$ws.Range("H2").Copy($ws.Range("I2"))
$ws.Range("H3").Copy($ws.Range("I3"))

# --- Row 5: add a status note in the new column ------------------------
$ws.Range("H2").Copy($ws.Range("I5"))
$ws.Range("I5").Value = "Introduction to RPA,  Introduction to UI Path and Installed UI path tool"

# --- Row 6: Shanth is on leave; GitHub doc update recorded -------------
$ws.Range("H6").Value = "On Leave"
$ws.Range("H6").WrapText = $true
$ws.Range("I6").Value = "GiTHub document is added and sent for Review. Working DatabaseTesting related document with Ramesh"
$ws.Range("I6").WrapText = $true

# --- Row heights recalculated for the new/updated wrapped content -----
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 90
$ws.Rows.Item(3).RowHeight = 90
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 90

# --- Restore the default top-left scroll position and move selection --
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("K6").Select()
